# Update the division-practice answer table.
# Each data row of the table (rows 1, 5, 9, 13, 17 -- the other rows are
# blank "work space" rows) holds 5 answers, one per column. We address the
# cells directly by (row, column) rather than using Find/Replace because
# several of the old/new strings collide with each other across cells
# (e.g. "79÷9=8, 7" is both an old value in one cell and a new value
# produced in another cell), which would make a blind global
# find-and-replace unsafe/ambiguous.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "86÷7=12, 2"
$t.Cell(1, 2).Range.Text  = "10÷7=1, 3"
$t.Cell(1, 3).Range.Text  = "90÷3=30, 0"
$t.Cell(1, 4).Range.Text  = "56÷8=7, 0"
$t.Cell(1, 5).Range.Text  = "78÷2=39, 0"

$t.Cell(5, 1).Range.Text  = "74÷3=24, 2"
$t.Cell(5, 2).Range.Text  = "17÷6=2, 5"
$t.Cell(5, 3).Range.Text  = "59÷7=8, 3"
$t.Cell(5, 4).Range.Text  = "56÷9=6, 2"
$t.Cell(5, 5).Range.Text  = "62÷8=7, 6"

$t.Cell(9, 1).Range.Text  = "99÷7=14, 1"
$t.Cell(9, 2).Range.Text  = "83÷6=13, 5"
$t.Cell(9, 3).Range.Text  = "40÷7=5, 5"
$t.Cell(9, 4).Range.Text  = "44÷6=7, 2"
$t.Cell(9, 5).Range.Text  = "79÷9=8, 7"

$t.Cell(13, 1).Range.Text = "75÷9=8, 3"
$t.Cell(13, 2).Range.Text = "71÷3=23, 2"
$t.Cell(13, 3).Range.Text = "77÷2=38, 1"
$t.Cell(13, 4).Range.Text = "72÷4=18, 0"
$t.Cell(13, 5).Range.Text = "89÷9=9, 8"

$t.Cell(17, 1).Range.Text = "59÷3=19, 2"
$t.Cell(17, 2).Range.Text = "33÷4=8, 1"
$t.Cell(17, 3).Range.Text = "29÷3=9, 2"
$t.Cell(17, 4).Range.Text = "59÷9=6, 5"
$t.Cell(17, 5).Range.Text = "25÷6=4, 1"

Write-Host "Updated 25 cells in table 1."
